$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 110.6
$ws.Range("F2").Value = 12.3
$ws.Range("G2").Value = 3.539333
$ws.Range("L2").Value = 16.5
$ws.Range("M2").Value = 0.976642

$ws.Range("E3").Value = 78.3
$ws.Range("F3").Value = 19.6
$ws.Range("G3").Value = 1.4815049
$ws.Range("H3").Value = 0.3356666666666666
$ws.Range("I3").Value = 0.326
$ws.Range("J3").Value = 0.3348333333333334
$ws.Range("K3").Value = 0.4351666666666667
$ws.Range("L3").Value = 12.16666666666667
$ws.Range("M3").Value = 1.42619

$ws.Range("E4").Value = 189.7
$ws.Range("G4").Value = 2.515083857142857
$ws.Range("H4").Value = 0.3464285714285715
$ws.Range("I4").Value = 0.3391428571428571
$ws.Range("J4").Value = 0.3291428571428572
$ws.Range("K4").Value = 0.4767142857142858
$ws.Range("L4").Value = 16.42857142857143
$ws.Range("M4").Value = 1.53299

$ws.Range("E5").Value = 77.2
$ws.Range("G5").Value = 2.05068
$ws.Range("H5").Value = 0.30175
$ws.Range("J5").Value = 0.314
$ws.Range("K5").Value = 0.36875
$ws.Range("L5").Value = 11.75
$ws.Range("M5").Value = 1.07535

$ws.Range("E6").Value = 54.3
$ws.Range("F6").Value = 18.1
$ws.Range("G6").Value = 1.8188089375
$ws.Range("H6").Value = 0.28375
$ws.Range("I6").Value = 0.29875
$ws.Range("J6").Value = 0.294125
$ws.Range("K6").Value = 0.341875
$ws.Range("L6").Value = 13.75
$ws.Range("M6").Value = 0.6144539999999999

$ws.Range("B7").Value = "young with upside"
$ws.Range("E7").Value = 79.59999999999999
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 1.702207066666667
$ws.Range("H7").Value = 0.31
$ws.Range("I7").Value = 0.315
$ws.Range("J7").Value = 0.3130000000000001
$ws.Range("K7").Value = 0.39325
$ws.Range("L7").Value = 15.5
$ws.Range("M7").Value = 2.90521

$ws.Range("E8").Value = 292.3
$ws.Range("F8").Value = 22.5
$ws.Range("G8").Value = 3.776856666666667
$ws.Range("H8").Value = 0.3283333333333333
$ws.Range("I8").Value = 0.3366666666666666
$ws.Range("J8").Value = 0.3233333333333333
$ws.Range("K8").Value = 0.4366666666666667
$ws.Range("L8").Value = 23.33333333333333
$ws.Range("M8").Value = 2.93718

$ws.Range("E9").Value = 171.2
$ws.Range("F9").Value = 17.1
$ws.Range("G9").Value = 3.8633445
$ws.Range("H9").Value = 0.347
$ws.Range("I9").Value = 0.35125
$ws.Range("J9").Value = 0.3465
$ws.Range("K9").Value = 0.45225
$ws.Range("L9").Value = 21.75
$ws.Range("M9").Value = 2.09657

$ws.Range("E10").Value = 302.9
$ws.Range("F10").Value = 33.7
$ws.Range("G10").Value = 2.525595
$ws.Range("H10").Value = 0.3755
$ws.Range("I10").Value = 0.369
$ws.Range("J10").Value = 0.3725
$ws.Range("K10").Value = 0.5044999999999999
$ws.Range("L10").Value = 25.5
$ws.Range("M10").Value = 3.84587

$ws.Range("E11").Value = 70.09999999999999
$ws.Range("F11").Value = 17.5
$ws.Range("G11").Value = 2.02550918
$ws.Range("H11").Value = 0.302
$ws.Range("I11").Value = 0.3128
$ws.Range("J11").Value = 0.301
$ws.Range("K11").Value = 0.3858
$ws.Range("L11").Value = 9.4
$ws.Range("M11").Value = 1.17187

$ws.Range("E12").Value = 222.7
$ws.Range("G12").Value = 3.2027375
$ws.Range("H12").Value = 0.3185
$ws.Range("I12").Value = 0.32525
$ws.Range("J12").Value = 0.31825
$ws.Range("K12").Value = 0.4062500000000001
$ws.Range("L12").Value = 16.75
$ws.Range("M12").Value = 4.02583

$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 97.8
$ws.Range("F13").Value = 24.5
$ws.Range("G13").Value = 0.8581548666666667
$ws.Range("H13").Value = 0.318
$ws.Range("I13").Value = 0.322
$ws.Range("J13").Value = 0.2985
$ws.Range("K13").Value = 0.437
$ws.Range("L13").Value = 22
$ws.Range("M13").Value = 1.84695

$ws.Range("E14").Value = 110.1
$ws.Range("F14").Value = 22
$ws.Range("G14").Value = 1.669761609090909
$ws.Range("H14").Value = 0.3719
$ws.Range("I14").Value = 0.3566
$ws.Range("J14").Value = 0.3591000000000001
$ws.Range("K14").Value = 0.5145000000000001
$ws.Range("M14").Value = 2.73367

$ws.Range("E15").Value = 180.5
$ws.Range("F15").Value = 22.6
$ws.Range("G15").Value = 3.13283325
$ws.Range("I15").Value = 0.351625
$ws.Range("J15").Value = 0.3575
$ws.Range("K15").Value = 0.50125
$ws.Range("M15").Value = 4.02683

$ws.Range("E16").Value = 72.90000000000001
$ws.Range("G16").Value = 1.645292333333333
$ws.Range("H16").Value = 0.2826666666666666
$ws.Range("I16").Value = 0.2776666666666667
$ws.Range("J16").Value = 0.2976666666666667
$ws.Range("K16").Value = 0.3396666666666667
$ws.Range("M16").Value = 2.19172

$ws.Range("E17").Value = 108.1
$ws.Range("G17").Value = 1.36938766
$ws.Range("H17").Value = 0.334
$ws.Range("I17").Value = 0.322
$ws.Range("J17").Value = 0.32
$ws.Range("K17").Value = 0.4503333333333333
$ws.Range("L17").Value = 22.33333333333333
$ws.Range("M17").Value = 1.57597

$ws.Range("E18").Value = 322.2
$ws.Range("F18").Value = 46
$ws.Range("G18").Value = 2.922371428571429
$ws.Range("H18").Value = 0.3681428571428572
$ws.Range("I18").Value = 0.3628571428571429
$ws.Range("J18").Value = 0.3491428571428571
$ws.Range("K18").Value = 0.5208571428571428
$ws.Range("M18").Value = 3.09539

$ws.Range("E19").Value = 284.3
$ws.Range("F19").Value = 21.9
$ws.Range("G19").Value = 2.571105
$ws.Range("H19").Value = 0.348
$ws.Range("I19").Value = 0.34325
$ws.Range("J19").Value = 0.33075
$ws.Range("K19").Value = 0.4742499999999999
$ws.Range("L19").Value = 20
$ws.Range("M19").Value = 2.8875

$ws.Range("E20").Value = 114.1
$ws.Range("F20").Value = 16.3
$ws.Range("G20").Value = 3.59123
$ws.Range("H20").Value = 0.3205
$ws.Range("I20").Value = 0.33425
$ws.Range("J20").Value = 0.3514999999999999
$ws.Range("K20").Value = 0.36625
$ws.Range("L20").Value = 9
$ws.Range("M20").Value = 2.26243

$ws.Range("B21").Value = "high-floor contributors"
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 64.5
$ws.Range("F21").Value = 12.9
$ws.Range("G21").Value = 2.31429
$ws.Range("H21").Value = 0.335
$ws.Range("I21").Value = 0.32
$ws.Range("K21").Value = 0.4385
$ws.Range("L21").Value = 20.5
$ws.Range("M21").Value = 1.6438
